$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell($table, $row, $col, $parts) {
    $cell = $table.Cell($row, $col)
    $innerXml = ""
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -gt 0) { $innerXml += "<w:br/>" }
        $innerXml += $parts[$i]
    }
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + $innerXml + '</w:r></w:p>'
    $cell.Range.InsertXML($xml)
}

Set-LatticeCell $t 1 1 @("<w:t>27 x 32</w:t>", "<w:t xml:space=`"preserve`">  3    2</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>2|    |</w:t>", "<w:t>7|    |</w:t>")
Set-LatticeCell $t 1 2 @("<w:t>13 x 40</w:t>", "<w:t xml:space=`"preserve`">  4    0</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>1|    |</w:t>", "<w:t>3|    |</w:t>")
Set-LatticeCell $t 1 3 @("<w:t>36 x 16</w:t>", "<w:t xml:space=`"preserve`">  1    6</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>3|    |</w:t>", "<w:t>6|    |</w:t>")
Set-LatticeCell $t 2 1 @("<w:t>86 x 67</w:t>", "<w:t xml:space=`"preserve`">  6    7</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>8|    |</w:t>", "<w:t>6|    |</w:t>")
Set-LatticeCell $t 2 2 @("<w:t>96 x 30</w:t>", "<w:t xml:space=`"preserve`">  3    0</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>9|    |</w:t>", "<w:t>6|    |</w:t>")
Set-LatticeCell $t 2 3 @("<w:t>71 x 61</w:t>", "<w:t xml:space=`"preserve`">  6    1</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>7|    |</w:t>", "<w:t>1|    |</w:t>")
Set-LatticeCell $t 3 1 @("<w:t>49 x 80</w:t>", "<w:t xml:space=`"preserve`">  8    0</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>4|    |</w:t>", "<w:t>9|    |</w:t>")
Set-LatticeCell $t 3 2 @("<w:t>18 x 88</w:t>", "<w:t xml:space=`"preserve`">  8    8</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>1|    |</w:t>", "<w:t>8|    |</w:t>")
Set-LatticeCell $t 3 3 @("<w:t>10 x 85</w:t>", "<w:t xml:space=`"preserve`">  8    5</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>1|    |</w:t>", "<w:t>0|    |</w:t>")
Set-LatticeCell $t 4 1 @("<w:t>26 x 16</w:t>", "<w:t xml:space=`"preserve`">  1    6</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>2|    |</w:t>", "<w:t>6|    |</w:t>")
Set-LatticeCell $t 4 2 @("<w:t>64 x 42</w:t>", "<w:t xml:space=`"preserve`">  4    2</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>6|    |</w:t>", "<w:t>4|    |</w:t>")
Set-LatticeCell $t 4 3 @("<w:t>66 x 81</w:t>", "<w:t xml:space=`"preserve`">  8    1</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>6|    |</w:t>", "<w:t>6|    |</w:t>")
Set-LatticeCell $t 5 1 @("<w:t>24 x 11</w:t>", "<w:t xml:space=`"preserve`">  1    1</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>2|    |</w:t>", "<w:t>4|    |</w:t>")
Set-LatticeCell $t 5 2 @("<w:t>85 x 43</w:t>", "<w:t xml:space=`"preserve`">  4    3</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>8|    |</w:t>", "<w:t>5|    |</w:t>")
Set-LatticeCell $t 5 3 @("<w:t>84 x 58</w:t>", "<w:t xml:space=`"preserve`">  5    8</w:t>", "<w:t xml:space=`"preserve`">  ----</w:t>", "<w:t>8|    |</w:t>", "<w:t>4|    |</w:t>")
